$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Data Isi Dunia"
